$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 5.5
$ws.Range("I5").Value = 1.48
$ws.Range("N5").Value = 1.4
$ws.Range("O5").Value = 2.88
$ws.Range("U5").Value = 34
$ws.Range("Y5").Value = 34
$ws.Range("AB5").Value = 17

# Row 15
$ws.Range("G15").Value = 2.25
$ws.Range("I15").Value = 3.25
$ws.Range("L15").Value = 1.44
$ws.Range("M15").Value = 2.63
$ws.Range("N15").Value = 2.4
$ws.Range("O15").Value = 1.53
$ws.Range("R15").Value = 2.05
$ws.Range("S15").Value = 1.7
$ws.Range("U15").Value = 9.5
$ws.Range("W15").Value = 21
$ws.Range("X15").Value = 21
$ws.Range("Z15").Value = 7
$ws.Range("AI15").Value = 34

# Row 16
$ws.Range("G16").Value = 2.15
$ws.Range("I16").Value = 3.3
$ws.Range("R16").Value = 1.91
$ws.Range("S16").Value = 1.91
$ws.Range("T16").Value = 7
$ws.Range("U16").Value = 10
$ws.Range("W16").Value = 21
$ws.Range("Y16").Value = 29
$ws.Range("Z16").Value = 9
$ws.Range("AB16").Value = 15
$ws.Range("AD16").Value = 301
$ws.Range("AE16").Value = 9.5
$ws.Range("AG16").Value = 12
$ws.Range("AH16").Value = 34

# Row 17
$ws.Range("J17").Value = 1.06
$ws.Range("K17").Value = 10
$ws.Range("L17").Value = 1.33
$ws.Range("M17").Value = 3.25

# Row 18
$ws.Range("K18").Value = 8

# Row 19
$ws.Range("L19").Value = 1.57
$ws.Range("M19").Value = 2.25
$ws.Range("N19").Value = 2.88
$ws.Range("O19").Value = 1.4
$ws.Range("P19").Value = 1.62
$ws.Range("Q19").Value = 2.2
$ws.Range("W19").Value = 26

# Row 23
$ws.Range("K23").Value = 7.5
$ws.Range("R23").Value = 2
$ws.Range("S23").Value = 1.73
$ws.Range("T23").Value = 6.5
$ws.Range("V23").Value = 10
$ws.Range("AC23").Value = 67
$ws.Range("AE23").Value = 8

# Row 35
$ws.Range("H35").Value = 3.65
$ws.Range("I35").Value = 4.75
$ws.Range("R35").Value = 1.82
$ws.Range("U35").Value = 7.6
$ws.Range("X35").Value = 13.5
$ws.Range("Y35").Value = 27
$ws.Range("AF35").Value = 27
$ws.Range("AG35").Value = 15.5

# Row 36
$ws.Range("G36").Value = 2.88
$ws.Range("I36").Value = 2.4

# Row 37
$ws.Range("G37").Value = 2.3
$ws.Range("I37").Value = 2.7
$ws.Range("T37").Value = 8.5
$ws.Range("W37").Value = 21
$ws.Range("X37").Value = 19
$ws.Range("AE37").Value = 9.5
$ws.Range("AF37").Value = 15
